$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Permutation mapping: new row number -> source row number (1-indexed data rows 2..12)
# Represents a row-shuffle of the dataset (weekly update), keeping columns
# A,B,C,E,F,G,H,I,N,O,Q,R constant and only moving D,J,K,L,M,P values.
$perm = @{
    2 = 3
    3 = 2
    4 = 8
    5 = 10
    6 = 6
    7 = 11
    8 = 5
    9 = 12
    10 = 7
    11 = 9
    12 = 4
}

# Snapshot original values for columns D, J, K, L, M, P before overwriting
$cols = @("D", "J", "K", "L", "M", "P")
$snapshot = @{}
foreach ($r in 2..12) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

foreach ($r in 2..12) {
    $src = $perm[$r]
    $srcData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcData[$c]
    }
}
